# edit.ps1 -- reproduce the target commit's edits via Word COM-interop.
#
# Several of the edits split existing sentences across multiple <w:r> runs,
# add a <w:proofErr> pair, or move a <w:bookmarkStart/End> pair -- structure
# that plain Find/Replace or Range.Text assignment cannot reproduce (this
# runtime coalesces same-formatted runs written that way back into one
# <w:r>). InsertXML, which accepts a WordOpenXML "single file package" and
# splices it into the target Range, preserves whatever run/bookmark/proofErr
# structure we hand it, so that's used throughout below.
#
# The document is processed from the BOTTOM up so that paragraph indices
# for content still to be visited never shift underneath us.

$d = $word.ActiveDocument

# U+2019 RIGHT SINGLE QUOTATION MARK ("wasn<’>t"). Built from a char code
# instead of being embedded literally, since this interpreter mangles
# non-ASCII literals placed directly in the script source.
$RSQUO = [char]8217

function New-PkgXml($innerBodyXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        $innerBodyXml + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function New-ListPara($ilvl, $runsXml) {
    return '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="' + $ilvl + `
        '"/><w:numId w:val="2"/></w:numPr></w:pPr>' + $runsXml + '</w:p>'
}

# =====================================================================
# 1) Last paragraph ("We could create a pie chart ...") is split in two:
#      a) a new paragraph about an additional table / bubble-ish idea
#         (contains a proofErr gramStart/gramEnd pair around "and  figure")
#      b) the original pie-chart sentence, reworded mid-sentence, with the
#         _GoBack bookmark now sitting later in the sentence.
#    Do this FIRST (it's the bottom of the document) so paragraph indices
#    used below are unaffected.
# ---------------------------------------------------------------------
$pPie = $d.Paragraphs(13)
# Make room: two new trailing paragraphs, inserted *before* either one gets
# its real content, so neither target Range is ever "the last paragraph in
# the document" when InsertXML runs on it (InsertXML on the document's
# current last paragraph otherwise leaves a stray empty trailing paragraph
# behind -- pre-splitting first avoids that).
$pPie.Range.InsertParagraphAfter()
$d.Paragraphs(14).Range.InsertParagraphAfter()

$additionalTableRuns = `
    '<w:r><w:t xml:space="preserve">An additional table we could create would include donation amount asked for in the Kickstarter and the genre of Kickstarter. These elements included, we could throw in the state of the campaign </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>and  figure</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> out whether there is a relationship between successful Kickstarter and lower asking amounts of funding.</w:t></w:r>'
$additionalTableXml = New-PkgXml(New-ListPara 1 $additionalTableRuns)
$d.Paragraphs(13).Range.InsertXML($additionalTableXml)

$pieChartRuns = `
    '<w:r><w:t xml:space="preserve">We could create a pie chart for total money funded among all Kickstarter campaigns based on category. This would </w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>help clearly visualize the popularity of specific classes of all Kickstarter campaigns.</w:t></w:r>'
$pieChartXml = New-PkgXml(New-ListPara 1 $pieChartRuns)
$d.Paragraphs(14).Range.InsertXML($pieChartXml)

# Drop the now-unneeded extra trailing paragraph created above.
$d.Paragraphs($d.Paragraphs.Count).Range.Delete()

# =====================================================================
# 2) "We could create tables ... campaigns and compare the results."
#    (which held the _GoBack bookmark) -> reworded as a single run, no
#    bookmark (the bookmark moved to the pie-chart paragraph above).
#    Must happen before anything that re-adds a _GoBack bookmark expects
#    id 0 to be free again, hence processed right after step 1.
# ---------------------------------------------------------------------
$bubbleRuns = '<w:r><w:t>We could create a bubble chart specifically for showcasing the difference in funding amounts for each parent category and subcategory.</w:t></w:r>'
$bubbleXml = New-PkgXml(New-ListPara 1 $bubbleRuns)
$d.Paragraphs(12).Range.InsertXML($bubbleXml)

# =====================================================================
# 3) "There wasn't information on social media impression ..." loses its
#    trailing clause and becomes two runs.
# ---------------------------------------------------------------------
$socialRuns = `
    '<w:r><w:t>There wasn' + $RSQUO + 't information on social media impression which could lead to further analysis/information</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">. </w:t></w:r>'
$socialXml = New-PkgXml(New-ListPara 1 $socialRuns)
$d.Paragraphs(10).Range.InsertXML($socialXml)

# =====================================================================
# 4) "Journalism were ..." -> "Journalism was ..." (typed as three runs),
#    plus a brand-new sub-bullet "Food trucks are very likely to fail"
#    inserted right after it, one indent level deeper.
# ---------------------------------------------------------------------
$pJournalism = $d.Paragraphs(7)
$pJournalism.Range.InsertParagraphAfter()

$foodTrucksRuns = '<w:r><w:t>Food trucks are very likely to fail</w:t></w:r>'
$foodTrucksXml = New-PkgXml(New-ListPara 2 $foodTrucksRuns)
$d.Paragraphs(8).Range.InsertXML($foodTrucksXml)

$journalismRuns = `
    '<w:r><w:t>Journalism w</w:t></w:r>' + `
    '<w:r><w:t>as</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> the least successful in gathering funds</w:t></w:r>'
$journalismXml = New-PkgXml(New-ListPara 1 $journalismRuns)
$d.Paragraphs(7).Range.InsertXML($journalismXml)

Write-Output ("Done. Final paragraph count: " + $d.Paragraphs.Count)
